$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old Q1 (DMN32D2LDF) row entirely; rows below shift up.
$ws.Rows.Item(21).Delete()

# Insert a new row at position 25 for the replacement Q1 part (QS5K2TR),
# pushing the LED row (and everything after) back down.
$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "QS5K2TR"
$ws.Range("C25").Value = "QS5K2TR"
$ws.Range("D25").Value = "SOT95P280X100-5N"
$ws.Range("E25").Value = "Q1"
$ws.Range("F25").Value = "2.5V Drive Nch+Nch MOSFET"

$ws.Range("A1").Select()
